$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Objeto" -> "Especificação" (label above the contract-object field)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Objeto", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Especificação", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Merge the two "4" runs that make up the "Saldo a Liquidar" amount
#    ("4" + "4" + ".500,00" -> "44" + ".500,00") while leaving the
#    ".500,00" run (which carries its own rsid) untouched/unmerged.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("44.500,00", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start

    # Temporarily drop Bold on the ".500,00" tail so its formatting no
    # longer matches the leading "4"+"4" runs - this stops the engine
    # from coalescing it together with them when we retype the "44".
    $tail = $d.Range($start + 2, $start + 9)
    $tail.Bold = 0

    # Re-type the first two characters ("4" + "4") as a single run "44".
    $d.Range($start, $start + 2).Select()
    $word.Selection.TypeText("44")

    # Restore the original (bold) formatting of the ".500,00" run - it
    # comes back as a clean "<w:b/>" with no left-over attributes.
    $tail2 = $d.Range($start + 2, $start + 9)
    $tail2.Bold = 1
}

# ---------------------------------------------------------------------
# 3) Mark the "Default Paragraph Font" style as hidden (adds
#    <w:semiHidden/> to its definition in styles.xml).
# ---------------------------------------------------------------------
$s = $d.Styles("Fontepargpadro")
try {
    $s.Hidden = $true
} catch {
    # Older/limited COM surfaces may not expose a settable Hidden
    # property on Style - ignore and continue.
}
